$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A (rows 2..97) currently holds labels "q1".."q96".
# Relabel them down by one: "q0".."q95".
for ($row = 2; $row -le 97; $row++) {
    $n = $row - 2
    $ws.Cells.Item($row, 1).Value = "q$n"
}
